$d = $word.ActiveDocument

function Replace-Exact($find, $replace) {
    $rng = $d.Content
    $rng.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null
}

# 1) The "Whack-A-Mole" title shape (DrawingML shape text)
$shape = $d.Shapes.Item(1)
$shape.TextFrame.TextRange.Text = "Whack-A-Prof"

# 2) Body paragraph replacements: "mole" -> "professor" (with matching casing/wording per
#    the original sentence), scoped tightly so each Find call targets a unique phrase.
Replace-Exact " and sound for our group “Whack-a-Mole” " " and sound for our group “Whack-a-Prof” "
Replace-Exact "“Whack-a-Mole” game implementations as a baseline " "“Whack-a-Prof” game implementations as a baseline "
Replace-Exact "and “moles”" "and “professors”"
Replace-Exact "what a “Whack-a-Mole” styled game is" "what a “Whack-a-Prof” styled game is"
Replace-Exact "can whack “moles”" "can whack “Professors”"
Replace-Exact "of moles is tied to the timer" "of Professors is tied to the timer"
Replace-Exact "as the session progresses “moles”" "as the session progresses “Professors”"
Replace-Exact "obtained per “mole” whacked" "obtained per “Professor” whacked"
Replace-Exact "“Bosses” are moles which occasionally " "“Bosses” are professors who occasionally "
Replace-Exact "including movement of “moles” " "including movement of “Professors” "
Replace-Exact "invokes a whack on the given “mole”. " "invokes a whack on the given “professor”. "
Replace-Exact "“Whack-a-Mole” project, it without a doubt will not be perfect " "“Whack-a-Prof” project, it without a doubt will not be perfect "
